# Update the "Förändrad" (changed) date column (C) for rows 2-11
# from 2023-09-01 (serial 45170) to 2023-09-05 (serial 45174).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45170) {
        $cell.Value2 = 45174
    }
}
